# "solved 4 ones problem"
# Fix up the "Preferred Date and Time" (column G) entries for a handful of
# members: two rows were missing a selection, two rows had ad-hoc free-text
# instead of one of the standard options, and one row had an incorrect
# leftover value that should be cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Person 4): was blank -> set to "Sunday Afternoon"
$ws.Range("G5").Value = "Sunday Afternoon"

# Row 10 (Person 9): was blank -> set to "Sunday Morning"
$ws.Range("G10").Value = "Sunday Morning"

# Row 17 (Person 16): had custom text "saturday afternoon" -> normalize to "Saturday Afternoon"
$ws.Range("G17").Value = "Saturday Afternoon"

# Row 18 (Person 17): had custom text "DRIVE SO TOLD HER SATURDAY MORNING ONLY" -> normalize to "Saturday Morning"
$ws.Range("G18").Value = "Saturday Morning"

# Row 20 (Person 19): incorrectly had "Saturday Morning" -> clear it out
$ws.Range("G20").ClearContents()
